$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume(1h) updates for existing rows (rates refreshed by the GitHub Action run) ---
$ws.Range("D2").Value = "37.120.29"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.047.79"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'247.89"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'56.26"
$ws.Range("E8").Value = "  -6.13%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "'0.0780"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'16.30"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "'0.885"
$ws.Range("E13").Value = "  +8.55%  "
$ws.Range("D14").Value = "2.346.13"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'5.70"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "2.053.74"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("E17").Value = "  +9.98%  "
$ws.Range("D18").Value = "37.124.77"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'74.44"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "'236.59"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "'169.85"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("D28").Value = "'20.08"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'4.93"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "'0.0618"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'4.47"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "'0.0882"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'2.24"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("D39").Value = "'5.25"
$ws.Range("E39").Value = "  +14.47%  "
$ws.Range("D40").Value = "'3.11"
$ws.Range("E40").Value = "  +8.79%  "
$ws.Range("D41").Value = "'0.0985"
$ws.Range("E41").Value = "  -17.03%  "
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "'95.59"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").Value = "1.267.76"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "'6.76"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "2.228.62"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").Value = "'44.03"
$ws.Range("E51").Value = "  -1.04%  "

# --- Rows 42/43 swap rank order: InjectiveProtocol now ranks above VeChain ---
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.53"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0222"
$ws.Range("E43").Value = "  -2.39%  "
